$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 4-9 (2015年-2020年) with more precise decimal values
$ws.Range("B4").Value = 38572.43
$ws.Range("C4").Value = 21446.16
$ws.Range("D4").Value = 29105.18
$ws.Range("E4").Value = 12230.85
$ws.Range("F4").Value = 65082.2

$ws.Range("B5").Value = 41805.58
$ws.Range("C5").Value = 23054.87
$ws.Range("D5").Value = 31521.77
$ws.Range("E5").Value = 13004.13
$ws.Range("F5").Value = 70347.78

$ws.Range("B6").Value = 45163.3978918216
$ws.Range("C6").Value = 24550.136696615
$ws.Range("D6").Value = 33781.3161775188
$ws.Range("E6").Value = 13723.0718494509
$ws.Range("F6").Value = 77097.1802596038

$ws.Range("B7").Value = 49173.5
$ws.Range("C7").Value = 24856.51
$ws.Range("D7").Value = 35196.11
$ws.Range("E7").Value = 14386.87
$ws.Range("F7").Value = 84907.13

$ws.Range("B8").Value = 52907.31
$ws.Range("C8").Value = 26783.67
$ws.Range("D8").Value = 37875.8
$ws.Range("E8").Value = 15549.37
$ws.Range("F8").Value = 91682.60000000001

$ws.Range("B9").Value = 54910.09
$ws.Range("C9").Value = 27501.14
$ws.Range("D9").Value = 39278.18
$ws.Range("E9").Value = 15597.71
$ws.Range("F9").Value = 96061.64

# Add new rows 10-11 (2021年, 2022年); copy the formatting (style) used by
# the preceding year-label cell (A9) into the new label cells A10:A11
$ws.Range("A9").Copy()
$ws.Range("A10:A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A10").Value = "2021年"
$ws.Range("B10").Value = 59005.2
$ws.Range("C10").Value = 30132.6
$ws.Range("D10").Value = 42498
$ws.Range("E10").Value = 16745.5
$ws.Range("F10").Value = 102595.8

$ws.Range("A11").Value = "2022年"
$ws.Range("B11").Value = 61724.13
$ws.Range("C11").Value = 31179.6
$ws.Range("D11").Value = 44282.85
$ws.Range("E11").Value = 16970.68
$ws.Range("F11").Value = 107224.07
